$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new value "***" into B14 and D14.
# B14 should inherit the same number-format style as the rest of column B (e.g. B2).
$ws.Range("B14").Value = "***"
$ws.Range("B14").NumberFormat = $ws.Range("B2").NumberFormat

$ws.Range("D14").Value = "***"

# Update the selected/active cell to D14 (matches the saved view state in the diff).
$ws.Range("D14").Select()
